$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the old "dbExcel" column
# (B -> C) and the old "WebExcel" column (C -> D) to the right, and carries
# their formatting/widths along automatically.
$ws.Columns("B").Insert()

# New column B header + query text (new test case: StatQuery).
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Female reproductive system cancer, NOS'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Match the wrapped-text style already used by the row-2 / column-A query cell.
$ws.Range("B2").WrapText = $true

# New column B should be as wide as column A (both hold long Cypher queries).
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Leave the selection where the author's edit ended up.
$null = $ws.Range("B7").Select()
